# Scheduled market-data refresh: updates currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) on a handful of Leve rows across several job
# sheets, reflecting newer Universalis price snapshots. Values only -
# no formulas, formatting, or structural changes involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce / Rubber
$ws.Range("H11").Value = 316.8
$ws.Range("I11").Value = 316.8
$ws.Range("K11").Value = 316.8
$ws.Range("M11").Value = -176.8

# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 1227.8889
$ws.Range("I18").Value = 1293
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1293
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -1009
$ws.Range("N18").Value = -1568

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 2883.5
$ws.Range("I76").Value = 2964.1177
$ws.Range("J76").Value = 2426.6667
$ws.Range("K76").Value = 2964.1177
$ws.Range("L76").Value = 2426.6667
$ws.Range("M76").Value = -2649.1177
$ws.Range("N76").Value = -3056.6667

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 2883.5
$ws.Range("I79").Value = 2964.1177
$ws.Range("J79").Value = 2426.6667
$ws.Range("K79").Value = 2964.1177
$ws.Range("L79").Value = 2426.6667
$ws.Range("M79").Value = -1872.1177
$ws.Range("N79").Value = -4610.6667

# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 651.6667
$ws.Range("I92").Value = 696.9231
$ws.Range("J92").Value = 534
$ws.Range("K92").Value = 696.9231
$ws.Range("L92").Value = 534
$ws.Range("M92").Value = 551.0769
$ws.Range("N92").Value = -3030

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 1765.2632
$ws.Range("I98").Value = 1208.2354
$ws.Range("K98").Value = 1208.2354
$ws.Range("M98").Value = 289.7646

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 1765.2632
$ws.Range("I122").Value = 1208.2354
$ws.Range("K122").Value = 3624.7062
$ws.Range("M122").Value = -1174.7062

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3962.1
$ws.Range("I137").Value = 3990
$ws.Range("J137").Value = 3870.4285
$ws.Range("K137").Value = 11970
$ws.Range("L137").Value = 11611.2855
$ws.Range("M137").Value = -9420
$ws.Range("N137").Value = -16711.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 919.8
$ws.Range("I2").Value = 810.8889
$ws.Range("J2").Value = 1900
$ws.Range("K2").Value = 810.8889
$ws.Range("L2").Value = 1900
$ws.Range("M2").Value = -697.8889
$ws.Range("N2").Value = -2126

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 20716.91
$ws.Range("I32").Value = 3752.3293
$ws.Range("J32").Value = 98000
$ws.Range("K32").Value = 3752.3293
$ws.Range("L32").Value = 98000
$ws.Range("M32").Value = -3465.3293
$ws.Range("N32").Value = -98574

# Row 39: Aurochs Star / Bull Hoplon
$ws.Range("H39").Value = 7132.5
$ws.Range("I39").Value = 7132.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 7132.5
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 843.8889
$ws.Range("I110").Value = 673.26086
$ws.Range("J110").Value = 1825
$ws.Range("K110").Value = 673.26086
$ws.Range("L110").Value = 1825
$ws.Range("M110").Value = 1371.73914
$ws.Range("N110").Value = -5915

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 919.8
$ws.Range("I116").Value = 810.8889
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 810.8889
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = 1483.1111
$ws.Range("N116").Value = -6488

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 919.8
$ws.Range("I3").Value = 810.8889
$ws.Range("J3").Value = 1900
$ws.Range("K3").Value = 810.8889
$ws.Range("L3").Value = 1900
$ws.Range("M3").Value = -696.8889
$ws.Range("N3").Value = -2128

# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 537.75
$ws.Range("I22").Value = 800.5
$ws.Range("J22").Value = 275
$ws.Range("K22").Value = 800.5
$ws.Range("L22").Value = 275
$ws.Range("M22").Value = -627.5
$ws.Range("N22").Value = -621

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 917.0476
$ws.Range("I94").Value = 902.1111
$ws.Range("J94").Value = 1006.6667
$ws.Range("K94").Value = 902.1111
$ws.Range("L94").Value = 1006.6667
$ws.Range("M94").Value = -451.1111
$ws.Range("N94").Value = -1908.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws.Range("H50").Value = 20538.182
$ws.Range("J50").Value = 20538.182
$ws.Range("L50").Value = 20538.182
$ws.Range("N50").Value = -21788.182

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 22.205883
$ws.Range("I2").Value = 10.217391
$ws.Range("J2").Value = 47.272728
$ws.Range("K2").Value = 61.304346
$ws.Range("L2").Value = 283.636368
$ws.Range("M2").Value = 51.695654
$ws.Range("N2").Value = -509.636368

# Row 26: A Grape Idea / Grape Juice
$ws.Range("H26").Value = 692.875
$ws.Range("I26").Value = 76.22221999999999
$ws.Range("J26").Value = 1485.7142
$ws.Range("K26").Value = 228.66666
$ws.Range("L26").Value = 4457.142599999999
$ws.Range("M26").Value = 59.33334000000002
$ws.Range("N26").Value = -5033.142599999999

# Row 38: Pretty as a Picture / Dark Vinegar
$ws.Range("H38").Value = 353.8
$ws.Range("I38").Value = 497
$ws.Range("J38").Value = 282.2
$ws.Range("K38").Value = 1491
$ws.Range("L38").Value = 846.5999999999999
$ws.Range("M38").Value = -1144
$ws.Range("N38").Value = -1540.6

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 1527.625
$ws.Range("I97").Value = 1462.6154
$ws.Range("J97").Value = 1604.4546
$ws.Range("K97").Value = 1462.6154
$ws.Range("L97").Value = 1604.4546
$ws.Range("M97").Value = -966.6153999999999
$ws.Range("N97").Value = -2596.4546

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 747.5714
$ws.Range("I22").Value = 594.5714
$ws.Range("J22").Value = 900.5714
$ws.Range("K22").Value = 594.5714
$ws.Range("L22").Value = 900.5714
$ws.Range("M22").Value = -299.5714
$ws.Range("N22").Value = -1490.5714

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 747.5714
$ws.Range("I27").Value = 594.5714
$ws.Range("J27").Value = 900.5714
$ws.Range("K27").Value = 594.5714
$ws.Range("L27").Value = 900.5714
$ws.Range("M27").Value = -487.5714
$ws.Range("N27").Value = -1114.5714

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 280504.28
$ws.Range("I46").Value = 2380.2
$ws.Range("J46").Value = 435017.66
$ws.Range("K46").Value = 2380.2
$ws.Range("L46").Value = 435017.66
$ws.Range("M46").Value = -2192.2
$ws.Range("N46").Value = -435393.66

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1188.9688
$ws.Range("I132").Value = 1074.1428
$ws.Range("J132").Value = 1221.12
$ws.Range("K132").Value = 3222.4284
$ws.Range("L132").Value = 3663.36
$ws.Range("M132").Value = -692.4284000000002
$ws.Range("N132").Value = -8723.360000000001
